# custom accuracy + new simulation run ("데이터 1000개"):
# refresh the 4 junction-flooding sample rows with new timestamps/values
# and drop the trailing row, shrinking the sheet from A1:AH6 to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (no longer present in the updated dataset; shifts dimension to A1:AH5)
$ws.Rows.Item(6).Delete()

# Replace data rows 2-5 (A:AH) with the refreshed simulation output
$row2 = @(45074.50694444445, 8.821999999999999, 6.243, 2.794, 19.605, 14.205, 5.893, 19.279, 10.957, 4.563, 5.697, 7.759, 8.683999999999999, 2.881, 7.143, 9.409000000000001, 6.779, 1.556, 0.752, 101.319, 19.6, 6.594, 12.192, 6.907, 0.798, 11.548, 5.824, 5.505, 6.107, 8.566000000000001, 2.095, 16.882, 3.332, 8.242000000000001)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @(45074.51388888889, 22.851, 17.001, 1.853, 50.152, 40.426, 17.384, 65.703, 27.839, 12.47, 17.767, 20.067, 21.483, 6.11, 18.046, 25.412, 15.504, 0.867, 0.962, 267.363, 50.403, 16.657, 33.573, 17.79, 2.318, 33.886, 14.713, 13.143, 15.358, 21.246, 0.973, 59.993, 9.24, 20.825)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @(45074.52083333334, 22.421, 16.745, 1.453, 49.122, 39.924, 17.223, 69.297, 27.267, 12.297, 17.662, 19.681, 20.999, 5.89, 17.67, 25.005, 15.032, 0.631, 0.865, 261.653, 49.44, 16.31, 33.095, 17.444, 2.293, 34.4, 14.407, 12.809, 15.017, 20.769, 0.644, 63.16, 9.116, 20.392)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @(45074.52777777778, 14.28, 10.66, 1, 31.36, 25.41, 10.92, 46.56, 17.38, 7.89, 11.18, 12.56, 13.43, 3.79, 11.28, 15.93, 9.65, 0.45, 0.58, 164.36, 31.62, 10.41, 21.11, 11.11, 1.46, 22.78, 9.199999999999999, 8.199999999999999, 9.609999999999999, 13.26, 0.48, 42.47, 5.79, 13.02)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# Update column widths that changed (stored XML width = ColumnWidth + 5/6, rounded to nearest 1/6)
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(18).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
